$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 155; this shifts all existing rows 155-302 down to 156-303,
# matching the source data (row302 -> row303, etc.) exactly as required by the diff.
$ws.Rows(155).Insert()

# Populate the newly inserted row 155 with its full record (columns A-R),
# reusing the constant values shared by every data row in this sheet.
$ws.Cells.Item(155, 1).Value = 5
$ws.Cells.Item(155, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(155, 3).Value = "Maule"
$ws.Cells.Item(155, 4).Value = 44669
$ws.Cells.Item(155, 5).Value = 7
$ws.Cells.Item(155, 6).Value = 100114014
$ws.Cells.Item(155, 7).Value = "Betarraga"
$ws.Cells.Item(155, 8).Value = "Sin especificar"
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 5000
$ws.Cells.Item(155, 11).Value = 600
$ws.Cells.Item(155, 12).Value = 600
$ws.Cells.Item(155, 13).Value = 600
$ws.Cells.Item(155, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(155, 15).Value = "Región del Maule"
$ws.Cells.Item(155, 16).Value = 120
$ws.Cells.Item(155, 17).Value = 5
$ws.Cells.Item(155, 18).Value = "Hortaliza"
